$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("total")

# ---------------------------------------------------------------------
# Add monthly subscription expenses (spotify, netflix) as new rows at
# the bottom of the "total" sheet (rows 1725-1746).
# ---------------------------------------------------------------------

$spotifyDates = @(45092,45122,45153,45184,45214,45245,45275,45306,45337,45366,45397)
$netflixDates = @(45082,45112,45143,45174,45204,45235,45265,45296,45327,45356,45387)

$firstRow = 1725

# --- Spotify: 11 monthly charges of 5.99, one per month Jun 2023 - Apr 2024 ---
for ($i = 0; $i -lt $spotifyDates.Count; $i++) {
    $r = $firstRow + $i

    $ws.Range("A$r").Value = $spotifyDates[$i]
    $ws.Range("B$r").Value = "entertainment"

    # store (J) was entered before the expense_type (C) for this transaction
    $ws.Range("J$r").Value = "spotify"
    $ws.Range("C$r").Value = "spotify subscription"

    $ws.Range("D$r").Value = 5.99

    $ws.Range("E$r").Formula = "=MONTH(A$r)"
    $ws.Range("F$r").Formula = "=YEAR(A$r)"
    $ws.Range("G$r").Formula = "=WEEKDAY(A$r, 2)"
    $ws.Range("H$r").Formula = '=CHOOSE(WEEKDAY(A' + $r + ', 2), "Monday", "Tuesday","Wednesday", "Thursday", "Friday", "Saturday","Sunday")'
    $ws.Range("I$r").Formula = '=TEXT(A' + $r + ', "MMM")'
}

$netflixFirstRow = $firstRow + $spotifyDates.Count

# --- Netflix: 11 monthly charges of 7.99, one per month Jun 2023 - Apr 2024 ---
for ($i = 0; $i -lt $netflixDates.Count; $i++) {
    $r = $netflixFirstRow + $i

    $ws.Range("A$r").Value = $netflixDates[$i]
    $ws.Range("B$r").Value = "entertainment"

    # expense_type (C) was entered before the store (J) for this transaction
    $ws.Range("C$r").Value = "netflix subscription"
    $ws.Range("J$r").Value = "netflix"

    $ws.Range("D$r").Value = 7.99

    $ws.Range("E$r").Formula = "=MONTH(A$r)"
    $ws.Range("F$r").Formula = "=YEAR(A$r)"
    $ws.Range("G$r").Formula = "=WEEKDAY(A$r, 2)"
    $ws.Range("H$r").Formula = '=CHOOSE(WEEKDAY(A' + $r + ', 2), "Monday", "Tuesday","Wednesday", "Thursday", "Friday", "Saturday","Sunday")'
    $ws.Range("I$r").Formula = '=TEXT(A' + $r + ', "MMM")'
}

$lastRow = $netflixFirstRow + $netflixDates.Count - 1

# Extend the (hidden) AutoFilter defined name along with the newly added data.
$wb.Names.Item("total!_FilterDatabase").RefersTo = "=total!`$A`$1:`$L`$1724"

# Reflect the final on-screen selection/state left after the edit.
$ws.Activate()
$ws.Range("J1737:J1746").Select()
